$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    (Heading1 "Play Book of Oz Free Today! Game Review") paragraph.
#    The paragraph is injected via a Flat-OPC InsertXML call so the run
#    layout (leading empty run + bold "Meta description" run + plain
#    rest-of-sentence run) matches exactly, rather than letting plain
#    text-insertion calls auto-merge adjacent same-formatted runs.
# ---------------------------------------------------------------------------
$metaParagraphXml = '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Explore the Wizard of Oz themed slot game, with respin and free spins features. Read the Book of Oz review and play for free today.</w:t></w:r></w:p>'

$documentXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $metaParagraphXml + '<w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document>'

$flatOpcXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + $documentXml + '</pkg:xmlData></pkg:part></pkg:package>'

$titlePara = $d.Paragraphs(1)
[void]$titlePara.Range.InsertXML($flatOpcXml, "After")

# ---------------------------------------------------------------------------
# 2) Remove the duplicate "Play Book of Oz Free Today! Game Review" paragraph
#    that used to sit near the end of the document (its content now lives in
#    the new Meta description paragraph up top). It is always the
#    second-to-last paragraph of the document, immediately before the
#    closing italic paragraph.
# ---------------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$dupeParagraph = $d.Paragraphs($paraCount - 1)
$dupeParagraph.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Update the remaining italic (last) paragraph's text -- previously the
#    meta-description sentence -- to the new "Prompt: ..." image-generation
#    prompt, keeping its existing italic run formatting untouched.
# ---------------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$lastParagraph = $d.Paragraphs($paraCount)
$promptRange = $d.Range($lastParagraph.Range.Start, $lastParagraph.Range.End)
$promptRange.Text = "Prompt: Create a colorful and eye-catching feature image for Book of Oz that showcases a happy Maya warrior with glasses in a cartoon style. The image should feature the emerald city in the background and the symbols of the game, such as the magic potions and the book symbol. Use a vibrant color scheme with green being the dominant color to reflect the theme of the game. Make sure the image stands out and catches the attention of potential players."
